$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that currently reads:
#   "Registration – adequate + Fb + email verification and Logout – Last seen"
# (split across two runs around a _GoBack bookmark) without assuming a
# fixed paragraph index.
# ------------------------------------------------------------------
$regIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pText = $d.Paragraphs.Item($i).Range.Text
    if ($pText -like "*Registration*adequate*Fb*email verification*") {
        $regIndex = $i
        break
    }
}

# ------------------------------------------------------------------
# Part 1: turn
#   [R1: "Registration – adequate + Fb + email verification "] bm [R2: "and Logout – Last seen"]
# into
#   [R1: "Registration "] [R2: "and Profile "] bm [R3: "– adequate + Fb + email verification and Logout – Last seen"]
# ------------------------------------------------------------------

# 1a. Shrink the first run down to just "Registration ".
$rngA = $d.Content
$oldFirst = "Registration – adequate + Fb + email verification "
$rngA.Find.Execute($oldFirst, $false, $false, $false, $false, $false, $true, 1, $false, "Registration ", 2)

# 1b. Insert a new run "and Profile " immediately before the bookmark. Splitting
# the paragraph in two at that point, typing into the (now separate) second
# paragraph and rejoining keeps "and Profile " as its own <w:r>, instead of
# being silently merged into the preceding run's text.
$bm = $d.Bookmarks.Item("_GoBack")
$splitPos = $bm.Start
$d.Range($splitPos, $splitPos).InsertParagraphAfter()

$tailPara = $d.Paragraphs.Item($regIndex + 1)
$tailPara.Range.InsertBefore("and Profile ")

$headPara = $d.Paragraphs.Item($regIndex)
$d.Range($headPara.Range.End - 1, $headPara.Range.End).Delete()

# 1c. Expand the final run's text to include the relocated "– adequate + Fb +
# email verification" prefix.
$rngB = $d.Content
$oldLast = "and Logout – Last seen"
$newLast = "– adequate + Fb + email verification and Logout – Last seen"
$rngB.Find.Execute($oldLast, $false, $false, $false, $false, $false, $true, 1, $false, $newLast, 2)

# ------------------------------------------------------------------
# Part 2: append a blank paragraph and three new paragraphs describing the
# Admin / course-filtering work after the Registration paragraph.
# ------------------------------------------------------------------

# 2a. Add a single trailing paragraph break after "Last seen" - using the
# Find/Replace "^p" token (rather than Range.InsertParagraphAfter) yields a
# genuinely empty <w:p/> instead of one containing a stray empty run.
$rngC = $d.Content
$rngC.Find.Execute($oldLast, $false, $false, $false, $false, $false, $true, 1, $false, $oldLast + "^p", 2)

$blankIndex = $regIndex + 1

# 2b. Build the three following paragraphs.
$pBlank = $d.Paragraphs.Item($blankIndex)
$pBlank.Range.InsertParagraphAfter()
$pAdmin = $d.Paragraphs.Item($blankIndex + 1)
$pAdmin.Range.InsertAfter("Admin - Filtering on the users and button functionality")
$pAdmin.Range.InsertParagraphAfter()

$pFilterCourses = $d.Paragraphs.Item($blankIndex + 2)
$pFilterCourses.Range.InsertAfter("Filtering the courses")
$pFilterCourses.Range.InsertParagraphAfter()

$pSort = $d.Paragraphs.Item($blankIndex + 3)
$pSort.Range.InsertAfter("Possibly sorting the courses and pagination")

# ------------------------------------------------------------------
# Part 3: split the "Admin - Filtering on the users and button functionality"
# paragraph into two runs ("Admin - " / "Filtering on the users and button
# functionality"), the same split + rejoin trick as Part 1b.
# ------------------------------------------------------------------
$pAdmin = $d.Paragraphs.Item($blankIndex + 1)
$prefix = "Admin - "
$splitPos2 = $pAdmin.Range.Start + $prefix.Length
$d.Range($splitPos2, $splitPos2).InsertParagraphAfter()

$pAdminHead = $d.Paragraphs.Item($blankIndex + 1)
$d.Range($pAdminHead.Range.End - 1, $pAdminHead.Range.End).Delete()

Write-Host "Done. Final text:"
Write-Host $d.Content.Text
